$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.976.88'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").Value = '3.540.79'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.38'
$ws.Range("E5").Value = '  -2.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.38'
$ws.Range("E6").Value = '  +4.44%  '
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.206'
$ws.Range("E9").Value = '  -5.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.653'
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.91'
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.52'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").Value = '4.102.38'
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '609.66'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.89'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D18").Value = '70.156.59'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").Value = '3.560.59'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.99'
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.31'
$ws.Range("E23").Value = '  +3.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.57'
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.62'
$ws.Range("E25").Value = '  -2.23%  '
$ws.Range("E26").Value = '  +2.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.62'
$ws.Range("E28").Value = '  -4.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.56'
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.34'
$ws.Range("E30").Value = '  +15.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.65'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("E33").Value = '  -2.15%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").Value = '0.0₃0858'
$ws.Range("E35").Value = '  +9.44%  '
$ws.Range("D36").Value = '3.744.65'
$ws.Range("E36").Value = '  +5.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.05'
$ws.Range("E38").Value = '  -3.78%  '
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.394'
$ws.Range("E40").Value = '  -1.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.67'
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '494.31'
$ws.Range("E42").Value = '  -8.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.134'
$ws.Range("E43").Value = '  -3.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0457'
$ws.Range("E44").Value = '  -2.73%  '
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("E46").Value = '  -4.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.65'
$ws.Range("E49").Value = '  -4.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000254'
$ws.Range("E50").Value = '  +4.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.60'
$ws.Range("E51").Value = '  -2.52%  '
